$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "261.57"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "1.58%"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "27.27"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "0.60%"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "4.726"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "2.64%"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.06067"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "2.84%"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.687"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "0.81%"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.8645"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "1.20%"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.9242"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "-2.28%"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1408"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "0.76%"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.05083"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "1.27%"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07162"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "1.08%"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.03035"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "-2.29%"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.09099"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.36%"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001540"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "0.80%"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0006081"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "0.79%"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "0.12%"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.449"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "-1.27%"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-0.43%"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "-1.26%"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.3127"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "2.38%"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "1.53%"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.096"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "3.66%"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04252"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "-0.37%"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "-0.39%"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.003911"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "-8.78%"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "0.03%"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "-18.94%"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.03882"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "1.54%"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1114"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "1.01%"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.004131"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "-33.58%"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.01491"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "4.93%"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.002204"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "-8.16%"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005345"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-0.07%"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-0.04%"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "6.93%"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "-47.51%"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "-0.04%"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.04%"

Write-Host "Applied updates to cryptos sheet"